# Auto-generated PowerShell COM-interop script
# Updates the cryptos list (Price / Volume(1h) columns, plus a row
# reorder for "Filecoin" / "ImmutableX") to match the refreshed
# GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        $Sheet,
        [string]$Address,
        [string]$Text
    )
    # Force the cell to Text format first so values that look like
    # numbers (e.g. "1.001", "49.70") are stored verbatim instead
    # of being parsed into a Double (which would silently drop
    # trailing zeros / significant digits). Then snap the style
    # back to "Normal" so no visible formatting/style is left
    # behind on the cell.
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextCell $ws 'D2' '22.391.06'
Set-TextCell $ws 'E2' '  +0.11%  '
Set-TextCell $ws 'D3' '1.570.88'
Set-TextCell $ws 'E3' '  +0.00%  '
Set-TextCell $ws 'E4' '  +0.01%  '
Set-TextCell $ws 'D5' '1.001'
Set-TextCell $ws 'E5' '  +0.04%  '
Set-TextCell $ws 'E6' '  +0.62%  '
Set-TextCell $ws 'D7' '0.3763'
Set-TextCell $ws 'E7' '  +2.23%  '
Set-TextCell $ws 'D8' '49.70'
Set-TextCell $ws 'E8' '  +0.38%  '
Set-TextCell $ws 'E9' '  +0.68%  '
Set-TextCell $ws 'D10' '0.07616'
Set-TextCell $ws 'E10' '  -0.04%  '
Set-TextCell $ws 'D11' '1.144'
Set-TextCell $ws 'E11' '  -1.71%  '
Set-TextCell $ws 'D12' '1.002'
Set-TextCell $ws 'E12' '  +0.01%  '
Set-TextCell $ws 'D13' '21.13'
Set-TextCell $ws 'E13' '  -0.86%  '
Set-TextCell $ws 'D14' '6.000'
Set-TextCell $ws 'E14' '  -1.00%  '
Set-TextCell $ws 'D15' '6.963'
Set-TextCell $ws 'D16' '1.569.30'
Set-TextCell $ws 'E16' '  -0.27%  '
Set-TextCell $ws 'E17' '  -0.07%  '
Set-TextCell $ws 'D18' '90.12'
Set-TextCell $ws 'E18' '  +0.63%  '
Set-TextCell $ws 'D19' '0.06744'
Set-TextCell $ws 'E19' '  +0.01%  '
Set-TextCell $ws 'E20' '  +0.02%  '
Set-TextCell $ws 'D21' '16.68'
Set-TextCell $ws 'E21' '  +1.05%  '
Set-TextCell $ws 'D22' '6.189'
Set-TextCell $ws 'E22' '  -0.66%  '
Set-TextCell $ws 'D23' '11.97'
Set-TextCell $ws 'E23' '  -0.17%  '
Set-TextCell $ws 'D24' '22.383.38'
Set-TextCell $ws 'E24' '  -0.05%  '
Set-TextCell $ws 'D25' '2.388'
Set-TextCell $ws 'E25' '  +0.44%  '
Set-TextCell $ws 'D26' '2.688'
Set-TextCell $ws 'E26' '  -7.53%  '
Set-TextCell $ws 'D27' '20.11'
Set-TextCell $ws 'E27' '  +0.51%  '
Set-TextCell $ws 'D28' '147.62'
Set-TextCell $ws 'E28' '  +1.01%  '
Set-TextCell $ws 'D29' '5.041'
Set-TextCell $ws 'E29' '  +1.55%  '
Set-TextCell $ws 'D30' '126.48'
Set-TextCell $ws 'E30' '  +0.67%  '
Set-TextCell $ws 'D31' '1.745.15'
Set-TextCell $ws 'E31' '  -0.06%  '
Set-TextCell $ws 'D32' '2.015'
Set-TextCell $ws 'E32' '  +0.65%  '
Set-TextCell $ws 'B33' 'ImmutableX'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D33' '0.9972'
Set-TextCell $ws 'E33' '  -4.42%  '
Set-TextCell $ws 'B34' 'Filecoin'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws 'D34' '6.079'
Set-TextCell $ws 'E34' '  -2.91%  '
Set-TextCell $ws 'D35' '10.08'
Set-TextCell $ws 'E35' '  -0.92%  '
Set-TextCell $ws 'D36' '0.08496'
Set-TextCell $ws 'E36' '  +0.43%  '
Set-TextCell $ws 'D37' '0.02532'
Set-TextCell $ws 'E37' '  -0.10%  '
Set-TextCell $ws 'D38' '1.382'
Set-TextCell $ws 'E38' '  +10.20%  '
Set-TextCell $ws 'D39' '0.2302'
Set-TextCell $ws 'E39' '  -0.98%  '
Set-TextCell $ws 'E40' '  -0.54%  '
Set-TextCell $ws 'D41' '5.405'
Set-TextCell $ws 'E41' '  -2.68%  '
Set-TextCell $ws 'D42' '11.38'
Set-TextCell $ws 'E42' '  -2.73%  '
Set-TextCell $ws 'D43' '0.6330'
Set-TextCell $ws 'D44' '1.001'
Set-TextCell $ws 'E44' '  +0.07%  '
Set-TextCell $ws 'D45' '13.94'
Set-TextCell $ws 'E45' '  -2.11%  '
Set-TextCell $ws 'D46' '3.804'
Set-TextCell $ws 'E46' '  +1.45%  '
Set-TextCell $ws 'D47' '0.5939'
Set-TextCell $ws 'E47' '  -0.70%  '
Set-TextCell $ws 'D48' '1.282'
Set-TextCell $ws 'E48' '  +1.45%  '
Set-TextCell $ws 'D49' '2.081'
Set-TextCell $ws 'E49' '  -1.43%  '
Set-TextCell $ws 'D50' '124.43'
Set-TextCell $ws 'E50' '  +0.24%  '
Set-TextCell $ws 'D51' '0.07319'
Set-TextCell $ws 'E51' '  +0.36%  '
